# Applies updated cryptocurrency price/volume data to Sheet1,
# matching the "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper that writes a value as plain text, preserving the exact
# literal characters (e.g. "1.000", "29.368.73") instead of letting
# Excel auto-convert number-looking text into a numeric value.
function Set-TextValue($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$ws.Cells.Item(2, 4).Value = "29.368.73"
$ws.Cells.Item(2, 5).Value = "  -0.14%  "
$ws.Cells.Item(3, 4).Value = "1.846.10"
$ws.Cells.Item(3, 5).Value = "  -0.21%  "
Set-TextValue 4 4 "0.9988"
$ws.Cells.Item(4, 5).Value = "  -0.19%  "
Set-TextValue 5 4 "240.68"
$ws.Cells.Item(5, 5).Value = "  +0.07%  "
Set-TextValue 6 4 "0.6306"
$ws.Cells.Item(6, 5).Value = "  +0.18%  "
Set-TextValue 7 4 "1.000"
$ws.Cells.Item(7, 5).Value = "  -0.03%  "
Set-TextValue 8 4 "0.07472"
$ws.Cells.Item(8, 5).Value = "  -2.24%  "
Set-TextValue 9 4 "0.2904"
$ws.Cells.Item(9, 5).Value = "  -0.25%  "
Set-TextValue 10 4 "24.38"
$ws.Cells.Item(10, 5).Value = "  -1.83%  "
Set-TextValue 11 4 "0.07730"
$ws.Cells.Item(11, 5).Value = "  -0.21%  "
$ws.Cells.Item(12, 4).Value = "1.846.14"
$ws.Cells.Item(12, 5).Value = "  -2.23%  "
Set-TextValue 13 4 "4.996"
$ws.Cells.Item(13, 5).Value = "  -0.74%  "
Set-TextValue 14 4 "0.6777"
$ws.Cells.Item(14, 5).Value = "  -0.50%  "
Set-TextValue 15 4 "0.00001020"
$ws.Cells.Item(15, 5).Value = "  -5.07%  "
Set-TextValue 16 4 "82.04"
$ws.Cells.Item(16, 5).Value = "  -1.54%  "
Set-TextValue 17 4 "6.134"
$ws.Cells.Item(17, 5).Value = "  -0.57%  "
$ws.Cells.Item(18, 4).Value = "29.429.44"
$ws.Cells.Item(18, 5).Value = "  -0.28%  "
$ws.Cells.Item(19, 5).Value = "  -0.15%  "
Set-TextValue 20 4 "12.29"
$ws.Cells.Item(20, 5).Value = "  -0.29%  "
Set-TextValue 21 4 "0.9999"
$ws.Cells.Item(21, 5).Value = "  -0.09%  "
Set-TextValue 22 4 "7.428"
$ws.Cells.Item(22, 5).Value = "  -0.54%  "
Set-TextValue 23 4 "1.002"
$ws.Cells.Item(23, 5).Value = "  +0.09%  "
Set-TextValue 24 4 "158.97"
$ws.Cells.Item(24, 5).Value = "  +0.61%  "
$ws.Cells.Item(25, 5).Value = "  -0.70%  "
Set-TextValue 26 4 "8.417"
$ws.Cells.Item(26, 5).Value = "  -0.14%  "
Set-TextValue 27 4 "17.53"
$ws.Cells.Item(27, 5).Value = "  -0.81%  "
Set-TextValue 28 4 "0.06354"
$ws.Cells.Item(28, 5).Value = "  +13.23%  "
Set-TextValue 29 4 "1.388"
$ws.Cells.Item(29, 5).Value = "  +0.86%  "
Set-TextValue 30 4 "1.474"
$ws.Cells.Item(30, 5).Value = "  +0.87%  "
$ws.Cells.Item(31, 5).Value = "  -1.00%  "
Set-TextValue 32 4 "4.057"
$ws.Cells.Item(32, 5).Value = "  +0.00%  "
$ws.Cells.Item(33, 5).Value = "  -1.21%  "
Set-TextValue 34 4 "1.141"
$ws.Cells.Item(34, 5).Value = "  -1.95%  "
Set-TextValue 35 4 "0.7003"
$ws.Cells.Item(35, 5).Value = "  +1.05%  "
Set-TextValue 36 4 "2.584"
$ws.Cells.Item(36, 5).Value = "  -0.20%  "
$ws.Cells.Item(37, 4).Value = "1.259.08"
$ws.Cells.Item(37, 5).Value = "  +2.56%  "
Set-TextValue 38 4 "2.836"
$ws.Cells.Item(38, 5).Value = "  +4.33%  "
Set-TextValue 39 4 "0.01817"
$ws.Cells.Item(39, 5).Value = "  +0.70%  "
Set-TextValue 40 4 "6.527"
$ws.Cells.Item(40, 5).Value = "  +1.21%  "
Set-TextValue 41 4 "0.9109"
$ws.Cells.Item(41, 5).Value = "  +0.47%  "
$ws.Cells.Item(42, 5).Value = "  -0.12%  "
$ws.Cells.Item(43, 4).Value = "2.007.69"
$ws.Cells.Item(43, 5).Value = "  -14.67%  "
Set-TextValue 44 4 "101.29"
$ws.Cells.Item(44, 5).Value = "  -0.44%  "
Set-TextValue 45 4 "66.29"
$ws.Cells.Item(45, 5).Value = "  +0.47%  "
Set-TextValue 46 4 "0.1171"
$ws.Cells.Item(46, 5).Value = "  +1.62%  "
Set-TextValue 47 4 "7.040"
$ws.Cells.Item(47, 5).Value = "  -2.11%  "
$ws.Cells.Item(48, 2).Value = "BabyDogeCoin"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue 48 4 "0.00000000116"
$ws.Cells.Item(48, 5).Value = "  -3.28%  "
$ws.Cells.Item(49, 2).Value = "EnergySwap"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue 49 4 "9.055"
$ws.Cells.Item(49, 5).Value = "  +0.74%  "
$ws.Cells.Item(50, 2).Value = "RenderToken"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue 50 4 "1.676"
$ws.Cells.Item(50, 5).Value = "  -0.34%  "
$ws.Cells.Item(51, 2).Value = "TheSandbox"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue 51 4 "0.3935"
$ws.Cells.Item(51, 5).Value = "  -2.17%  "
